$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quarter header / publish-date labels in row 9 (I, J, M)
$ws.Range("I9").Value = "1402-01-28 (5)"
$ws.Range("J9").Value = "1402-01-28 (8)"
$ws.Range("M9").Value = "1402-01-28 (3)"

# Update EPS-after-tax row (row 25): previously "-" placeholders, now computed values
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("M25").Value = 0

# Update Capital row (row 26): previously "-" placeholders, now computed values
$ws.Range("I26").Value = 7688
$ws.Range("J26").Value = 7580
$ws.Range("M26").Value = 6289
